# "Add files via upload" – renumbers Proposiciones I.1-I.16 to Roman numerals
# (I, II, III, ... XVI), tweaks three of their descriptions, and moves the
# active selection from B18 to B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Proposición I.1..I.16 -> Proposición I.I..I.XVI -----------
$romanLabels = @{
    2  = "Proposición I.I"
    3  = "Proposición I.II"
    4  = "Proposición I.III"
    5  = "Proposición I.IV"
    6  = "Proposición I.V"
    7  = "Proposición I.VI"
    8  = "Proposición I.VII"
    9  = "Proposición I.VIII"
    10 = "Proposición I.IX"
    11 = "Proposición I.X"
    12 = "Proposición I.XI"
    13 = "Proposición I.XII"
    14 = "Proposición I.XIII"
    15 = "Proposición I.XIV"
    16 = "Proposición I.XV"
    17 = "Proposición I.XVI"
}

foreach ($row in $romanLabels.Keys) {
    $ws.Cells.Item($row, 1).Value = $romanLabels[$row]
}

# --- Column B: description text updates for rows 9, 13 and 17 -----------
$ws.Range("B9").Value = "Si dos triángulos tienen dos lados del uno iguales respectivamente a dos lados del otro y tienen también iguales sus bases respectivas, también tendrán iguales los ángulos comprendidos entre las rectas."
$ws.Range("B13").Value = "Trazar una línea recta perpendicular a una recta infinita dada desde un punto que no está en ella."
$ws.Range("B17").Value = "Para cualquier triángulo, cuando un lado se prolonga, el ángulo externo es mayor que cada uno de los ángulos internos opuestos."

# --- Row 9 now wraps onto a second line, so its height grows ------------
$ws.Rows.Item(9).RowHeight = 28.8

# --- Column A widens slightly to fit the new longest label ("Proposición I.VIII") ---
$ws.Columns.Item(1).ColumnWidth = 14

# --- Selection / view: move from B18 to B17, drop the frozen scroll -----
$ws.Range("B17").Select() | Out-Null
